# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
#
# The source feed had four pairs of match-rows whose data got crossed
# (results / odds recorded against the wrong fixture row). This swaps the
# full row content (every column from B to AB - i.e. everything except the
# running "id" in column A) between each pair of rows, restoring the
# correct pairing of HomeTeam/AwayTeam/result/odds per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($r1, $r2) {
    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")
    $temp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $temp
}

Swap-RowData 16 17
Swap-RowData 67 68
Swap-RowData 86 87
Swap-RowData 161 162
